# ============================================================================
# feat: add 2022-Q3 data
#
# 1. Insert a new worksheet "2022-Q3" right after "总计" (i.e. before the
#    existing "2022-Q2" sheet) and populate it with the quarter's fund
#    holdings table.
# 2. Update the "总计" (summary) sheet: insert a new top row for 2022-Q3 and
#    shift the existing quarterly rows down by one.
# ============================================================================

$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------------
# Step 1: create the new "2022-Q3" worksheet right after "总计"
# ----------------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q3Sheet.Name = "2022-Q3"

# Header row (B1:H1), matching the style used on every other quarter sheet
$q3Headers = @(
    "基金代码",
    "基金名称",
    "基金规模",
    "股票总仓位",
    "仓位占比",
    "持有市值(亿元)",
    "仓位排名"
)
for ($i = 0; $i -lt $q3Headers.Count; $i++) {
    $col = 2 + $i   # column B == 2
    $q3Sheet.Cells.Item(1, $col).Value = $q3Headers[$i]
}
# Reuse the existing bold/border/center header style from the summary sheet
$totalSheet.Range("B1").Copy()
$q3Sheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Data rows 2-8: A(index,n) B(code,text) C(name,text) D(size,text)
# E(position,text) F(ratio,text) G(value,text) H(rank,n)
$q3Rows = @(
    ,@(0, "257010", "国联安小盘精选混合", "8.50", "74.70", "5.56", "0.4726", 4)
    ,@(1, "006138", "国联安价值优选股票", "0.57", "94.64", "5.35", "0.0305", 7)
    ,@(2, "011243", "万家惠裕回报6个月持有期混合A", "1.54", "27.67", "1.39", "0.0214", 1)
    ,@(3, "002367", "国联安安稳灵活配置混合", "0.57", "47.79", "3.00", "0.0171", 4)
    ,@(4, "011244", "万家惠裕回报6个月持有期混合C", "0.12", "27.67", "1.39", "0.0017", 1)
    ,@(5, "004791", "富荣中证500指数增强C", "0.08", "89.53", "1.78", "0.0014", 10)
    ,@(6, "004790", "富荣中证500指数增强A", "0.02", "89.53", "1.78", "0.0004", 10)
)

# Columns B,D,E,F,G hold numeric-looking text (fund codes / formatted
# percentages) that must stay literal text (leading zeros, trailing zeros),
# so format those columns as Text before writing the values.
$q3Sheet.Range("B2:B8").NumberFormat = "@"
$q3Sheet.Range("D2:G8").NumberFormat = "@"

for ($r = 0; $r -lt $q3Rows.Count; $r++) {
    $row = $q3Rows[$r]
    $excelRow = $r + 2
    $q3Sheet.Cells.Item($excelRow, 1).Value = $row[0]   # A - numeric index
    $q3Sheet.Cells.Item($excelRow, 2).Value = $row[1]   # B - fund code (text)
    $q3Sheet.Cells.Item($excelRow, 3).Value = $row[2]   # C - fund name
    $q3Sheet.Cells.Item($excelRow, 4).Value = $row[3]   # D - fund size (text)
    $q3Sheet.Cells.Item($excelRow, 5).Value = $row[4]   # E - stock position (text)
    $q3Sheet.Cells.Item($excelRow, 6).Value = $row[5]   # F - position ratio (text)
    $q3Sheet.Cells.Item($excelRow, 7).Value = $row[6]   # G - held value (text)
    $q3Sheet.Cells.Item($excelRow, 8).Value = $row[7]   # H - position rank (n)
}

# Column A (row index) carries the same bold/border/center style as the
# other quarter sheets.
$totalSheet.Range("A2").Copy()
$q3Sheet.Range("A2:A8").PasteSpecial($xlPasteFormats)

# ----------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet - insert the 2022-Q3 row at the top
# of the data (row 2) and shift every other quarter down by one row.
# ----------------------------------------------------------------------------
# Final state for rows 2-9: A(index,n) B(quarter,text) C(count,n) D(value,n)
$totalRows = @(
    ,@(0, "2022-Q3", 7, 0.55)
    ,@(1, "2022-Q2", 12, 1.61)
    ,@(2, "2022-Q1", 6, 0.73)
    ,@(3, "2021-Q4", 7, 1.08)
    ,@(4, "2021-Q3", 5, 1.17)
    ,@(5, "2021-Q2", 12, 7.04)
    ,@(6, "2021-Q1", 19, 11.75)
    ,@(7, "2020-Q4", 13, 3.29)
)

# Give new row 9 (beyond the old A1:D8 dimension) the same column-A style
# used by the rest of the table before writing any values into it.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A9").PasteSpecial($xlPasteFormats)

for ($r = 0; $r -lt $totalRows.Count; $r++) {
    $row = $totalRows[$r]
    $excelRow = $r + 2
    $totalSheet.Cells.Item($excelRow, 1).Value = $row[0]   # A - index
    $totalSheet.Cells.Item($excelRow, 2).Value = $row[1]   # B - quarter label
    $totalSheet.Cells.Item($excelRow, 3).Value = $row[2]   # C - holding count
    $totalSheet.Cells.Item($excelRow, 4).Value = $row[3]   # D - holding value
}
